# Update XLS with Code School course completions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Rails for Zombies 2" - mark date done (9/28/2013 -> serial 41545)
$ws.Range("E39").Value = 41545

# "Rails Best Practices" - mark date done (9/29/2013 -> serial 41546)
$ws.Range("E42").Value = 41546

# "Rails Testing with Rspec" - mark date done (9/29/2013 -> serial 41546)
$ws.Range("E43").Value = 41546

# Move the active selection to E30, matching where the author left off editing.
$ws.Range("E30").Select()
